$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35 - this shifts the existing rows 35..142
# down to 36..143 (matching the target dimension A1:T143).
$ws.Rows(35).Insert()

# Populate the newly inserted row 35 with the new record's data.
$ws.Range("A35").Value = 4
$ws.Range("B35").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C35").Value = "Los Lagos"
$ws.Range("D35").Value = 45037
$ws.Range("E35").Value = 10
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100104
$ws.Range("H35").Value = "Frutos de pepita"
$ws.Range("I35").Value = 100104003
$ws.Range("J35").Value = "Membrillo"
$ws.Range("K35").Value = "Champion"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 400
$ws.Range("N35").Value = 14000
$ws.Range("O35").Value = 15000
$ws.Range("P35").Value = 14500
$ws.Range("Q35").Value = "$/caja 18 kilos empedrada"
$ws.Range("R35").Value = "Región de O'Higgins"
$ws.Range("S35").Value = 806
$ws.Range("T35").Value = 18
